$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a value to be stored as text, even if it looks like a number,
# without leaving any residual cell-style change behind.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row -> Price (D) / Volume(1h) (E) updates, matching the authoritative diff.
$updates = @{
    2  = @{ D = "26.277.02"; E = "  +0.54%  " }
    3  = @{ D = "1.663.65";  E = "  +0.46%  " }
    4  = @{ E = "  +0.87%  " }
    5  = @{ D = "218.78";    E = "  +0.21%  " }
    6  = @{ D = "0.5308";    E = "  +0.30%  " }
    7  = @{ E = "  +0.81%  " }
    8  = @{ D = "0.2635";    E = "  +1.03%  " }
    9  = @{ D = "0.06362";   E = "  +0.37%  " }
    10 = @{ D = "20.54";     E = "  +0.56%  " }
    11 = @{ D = "0.07850";   E = "  +0.97%  " }
    12 = @{ D = "4.563";     E = "  +1.48%  " }
    13 = @{ D = "1.667.78";  E = "  +1.01%  " }
    14 = @{ D = "1.892.53";  E = "  +0.49%  " }
    15 = @{ D = "0.5527" }
    16 = @{ D = "0.0₅8178";  E = "  +0.09%  " }
    17 = @{ E = "  +0.32%  " }
    19 = @{ E = "  +2.43%  " }
    20 = @{ D = "192.77";    E = "  -0.02%  " }
    21 = @{ E = "  +1.28%  " }
    22 = @{ D = "6.031";     E = "  +0.07%  " }
    23 = @{ E = "  +0.85%  " }
    24 = @{ D = "144.36";    E = "  +1.71%  " }
    25 = @{ E = "  -2.03%  " }
    26 = @{ D = "7.199";     E = "  -1.01%  " }
    27 = @{ D = "16.07";     E = "  -0.80%  " }
    28 = @{ E = "  +2.79%  " }
    29 = @{ D = "0.05924";   E = "  -0.22%  " }
    30 = @{ D = "1.282";     E = "  +0.16%  " }
    31 = @{ D = "3.594";     E = "  +2.17%  " }
    32 = @{ E = "  +0.97%  " }
    33 = @{ D = "1.615";     E = "  +2.31%  " }
    34 = @{ D = "2.828";     E = "  +1.22%  " }
    35 = @{ D = "0.9595";    E = "  +1.07%  " }
    36 = @{ D = "2.426";     E = "  +0.65%  " }
    37 = @{ D = "0.5795";    E = "  +2.25%  " }
    38 = @{ D = "0.01603";   E = "  -0.46%  " }
    39 = @{ D = "0.8660";    E = "  +2.16%  " }
    40 = @{ D = "5.860";     E = "  +0.80%  " }
    41 = @{ E = "  +0.78%  " }
    42 = @{ D = "1.046.64";  E = "  +2.30%  " }
    43 = @{ D = "104.08";    E = "  +1.40%  " }
    44 = @{ D = "1.805.47";  E = "  +0.37%  " }
    45 = @{ E = "  +0.38%  " }
    46 = @{ E = "  -5.14%  " }
    47 = @{ E = "  +0.64%  " }
    48 = @{ E = "  +2.22%  " }
    49 = @{ D = "7.977";     E = "  +2.68%  " }
    50 = @{ D = "0.05166";   E = "  +0.31%  " }
    51 = @{ D = "1.434";     E = "  -3.00%  " }
}

foreach ($rowNum in $updates.Keys) {
    $rowData = $updates[$rowNum]
    if ($rowData.ContainsKey("D")) {
        Set-TextValue $ws.Range("D$rowNum") $rowData["D"]
    }
    if ($rowData.ContainsKey("E")) {
        $ws.Range("E$rowNum").Value = $rowData["E"]
    }
}
